# release verze 4.0.0, oprava dopadů zoomu na ip setting
#
# - switch the active/selected sheet from "Settings_recources" back to "Settings"
# - Settings: update a few numeric option cells
# - Settings_recources: restore the previous column width, move the sheet's
#   selection back to A1, and refresh several resource/config values
#   (file type lists now stored as python-list literals, updated user paths,
#   zoom settings, rendering mode, ...)

$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("Settings")
$wsResources = $wb.Worksheets.Item("Settings_recources")

# ---- Settings sheet -------------------------------------------------
$wsSettings.Range("B1").Value = 4
$wsSettings.Range("B3").Value = 1
$wsSettings.Range("B4").Value = 0
$wsSettings.Range("B5").Value = 0

# ---- Settings_recources sheet ---------------------------------------
$wsResources.Range("B1").Value = "['bmp', 'png']"
$wsResources.Range("B2").Value = "['jpg', 'bmp', 'png', 'ifz']"
$wsResources.Range("B3").Value = "C:/Users/jakub.hlavacek.local/Pictures/Screenshots/"
# these four/six cells hold numbers-as-text in the source sheet (no numeric
# cell type) - a leading apostrophe keeps Excel from re-typing them as
# real numbers when assigned through automation
$wsResources.Range("B4").Value = "'998"
$wsResources.Range("B5").Value = "['28', '02', '2024']"
$wsResources.Range("B9").Value = "'55"
$wsResources.Range("B19").Value = "'30"
$wsResources.Range("B20").Value = "'85"
$wsResources.Range("B23").Value = "'8"
$wsResources.Range("B29").Value = "C:/Users/jakub.hlavacek.local/Desktop/JHV/Work/TRIMAZKON/"
$wsResources.Range("A30").Value = "nastavení celkového zoomu [%]:"
$wsResources.Range("B30").Value = "'60"
$wsResources.Range("B31").Value = "ano"
$wsResources.Range("A32").Value = "Nastavení módu vykreslování (Katalog)"
$wsResources.Range("B32").Value = "fast"

# restore the default column width (drop the custom 71.85... width) on column A
$wsResources.Columns.Item(1).ColumnWidth = 8.43

# ---- window/selection state ------------------------------------------
# move Settings_recources' own cursor back to A1
$wsResources.Activate()
$wsResources.Range("A1").Select()

# Settings becomes the active/selected sheet again (activeTab 5 -> 3),
# keeping its previous selection at B24
$wsSettings.Activate()
$wsSettings.Range("B24").Select()
